# Scheduled runner: refresh Leve profit calculations with updated market prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1322.4546
$ws.Range("J112").Value = 1322.4546
$ws.Range("L112").Value = 3967.3638
$ws.Range("N112").Value = -6183.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6050.087
$ws.Range("I32").Value = 4428.8086
$ws.Range("K32").Value = 4428.8086
$ws.Range("M32").Value = -4141.8086
$ws.Range("H45").Value = 4335
$ws.Range("I45").Value = 5505.5
$ws.Range("J45").Value = 1994
$ws.Range("K45").Value = 5505.5
$ws.Range("L45").Value = 1994
$ws.Range("M45").Value = -5128.5
$ws.Range("N45").Value = -2748
$ws.Range("H52").Value = 25440
$ws.Range("J52").Value = 25440
$ws.Range("L52").Value = 25440
$ws.Range("N52").Value = -26076
$ws.Range("H103").Value = 34934.24
$ws.Range("J103").Value = 34934.24
$ws.Range("L103").Value = 34934.24
$ws.Range("N103").Value = -37278.24
$ws.Range("H122").Value = 8610.444
$ws.Range("I122").Value = 7246.6665
$ws.Range("J122").Value = 11338
$ws.Range("K122").Value = 21739.9995
$ws.Range("L122").Value = 34014
$ws.Range("M122").Value = -19289.9995
$ws.Range("N122").Value = -38914
$ws.Range("H132").Value = 2505.524
$ws.Range("I132").Value = 1174.1666
$ws.Range("J132").Value = 4280.6665
$ws.Range("K132").Value = 3522.4998
$ws.Range("L132").Value = 12841.9995
$ws.Range("M132").Value = -992.4998000000001
$ws.Range("N132").Value = -17901.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2365.8462
$ws.Range("I134").Value = 1385.6364
$ws.Range("J134").Value = 7757
$ws.Range("K134").Value = 4156.9092
$ws.Range("L134").Value = 23271
$ws.Range("M134").Value = -1621.9092
$ws.Range("N134").Value = -28341

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7681.5
$ws.Range("I31").Value = 1562.2222
$ws.Range("K31").Value = 1562.2222
$ws.Range("M31").Value = -1267.2222
$ws.Range("H34").Value = 7681.5
$ws.Range("I34").Value = 1562.2222
$ws.Range("K34").Value = 1562.2222
$ws.Range("M34").Value = -1360.2222
$ws.Range("H44").Value = 20000
$ws.Range("J44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -20884
$ws.Range("H48").Value = 39167
$ws.Range("J48").Value = 39167
$ws.Range("L48").Value = 39167
$ws.Range("N48").Value = -40119
$ws.Range("H50").Value = 29477.273
$ws.Range("J50").Value = 29477.273
$ws.Range("L50").Value = 29477.273
$ws.Range("N50").Value = -30727.273
$ws.Range("H59").Value = 35019
$ws.Range("J59").Value = 35019
$ws.Range("L59").Value = 35019
$ws.Range("N59").Value = -37309
$ws.Range("H60").Value = 15604.448
$ws.Range("J60").Value = 15604.448
$ws.Range("L60").Value = 15604.448
$ws.Range("N60").Value = -16626.448
$ws.Range("H68").Value = 62625.5
$ws.Range("J68").Value = 62625.5
$ws.Range("L68").Value = 62625.5
$ws.Range("N68").Value = -64123.5
$ws.Range("H71").Value = 62625.5
$ws.Range("J71").Value = 62625.5
$ws.Range("L71").Value = 187876.5
$ws.Range("N71").Value = -195364.5
$ws.Range("H86").Value = 4149.8335
$ws.Range("I86").Value = 4450
$ws.Range("J86").Value = 3999.75
$ws.Range("K86").Value = 4450
$ws.Range("L86").Value = 3999.75
$ws.Range("M86").Value = -3327
$ws.Range("N86").Value = -6245.75
$ws.Range("H89").Value = 4149.8335
$ws.Range("I89").Value = 4450
$ws.Range("J89").Value = 3999.75
$ws.Range("K89").Value = 22250
$ws.Range("L89").Value = 19998.75
$ws.Range("M89").Value = -16634
$ws.Range("N89").Value = -31230.75
$ws.Range("H132").Value = 2822.3157
$ws.Range("I132").Value = 1708.4
$ws.Range("K132").Value = 5125.200000000001
$ws.Range("M132").Value = -2595.200000000001
$ws.Range("H137").Value = 33239.855
$ws.Range("J137").Value = 33239.855
$ws.Range("L137").Value = 33239.855
$ws.Range("N137").Value = -43439.855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1900.9286
$ws.Range("I64").Value = 1333
$ws.Range("J64").Value = 2055.818
$ws.Range("K64").Value = 3999
$ws.Range("L64").Value = 6167.454000000001
$ws.Range("M64").Value = -3729
$ws.Range("N64").Value = -6707.454000000001
$ws.Range("H67").Value = 1900.9286
$ws.Range("I67").Value = 1333
$ws.Range("J67").Value = 2055.818
$ws.Range("K67").Value = 3999
$ws.Range("L67").Value = 6167.454000000001
$ws.Range("M67").Value = -3063
$ws.Range("N67").Value = -8039.454000000001
$ws.Range("H95").Value = 5000
$ws.Range("J95").Value = 5000
$ws.Range("L95").Value = 15000
$ws.Range("N95").Value = -19118
$ws.Range("H113").Value = 4032841
$ws.Range("I113").Value = 592.5238000000001
$ws.Range("J113").Value = 12500563
$ws.Range("K113").Value = 1777.5714
$ws.Range("L113").Value = 37501689
$ws.Range("M113").Value = 392.4285999999997
$ws.Range("N113").Value = -37506029
$ws.Range("H131").Value = 779.05
$ws.Range("J131").Value = 830.3333
$ws.Range("L131").Value = 2490.9999
$ws.Range("N131").Value = -12570.9999
$ws.Range("H141").Value = 8712.23
$ws.Range("I141").Value = 9585.9
$ws.Range("J141").Value = 5800
$ws.Range("K141").Value = 28757.7
$ws.Range("L141").Value = 17400
$ws.Range("M141").Value = -23577.7
$ws.Range("N141").Value = -27760

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5764.375
$ws.Range("I70").Value = 5428.0464
$ws.Range("J70").Value = 6876.846
$ws.Range("K70").Value = 5428.0464
$ws.Range("L70").Value = 6876.846
$ws.Range("M70").Value = -5158.0464
$ws.Range("N70").Value = -7416.846
$ws.Range("H73").Value = 5764.375
$ws.Range("I73").Value = 5428.0464
$ws.Range("J73").Value = 6876.846
$ws.Range("K73").Value = 5428.0464
$ws.Range("L73").Value = 6876.846
$ws.Range("M73").Value = -4492.0464
$ws.Range("N73").Value = -8748.846
$ws.Range("H97").Value = 2010
$ws.Range("I97").Value = 2010
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 2010
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1514
$ws.Range("N97").ClearContents()
$ws.Range("H102").Value = 3157.625
$ws.Range("I102").Value = 2202.2
$ws.Range("K102").Value = 2202.2
$ws.Range("M102").Value = -580.1999999999998
$ws.Range("H113").Value = 1633.3334
$ws.Range("I113").Value = 1633.3334
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1633.3334
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 536.6666
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 3450
$ws.Range("I126").Value = 2876.7124
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 8630.137199999999
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -6160.137199999999
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9750
$ws.Range("I40").Value = 9000
$ws.Range("J40").Value = 10500
$ws.Range("K40").Value = 9000
$ws.Range("L40").Value = 10500
$ws.Range("M40").Value = -8864
$ws.Range("N40").Value = -10772
$ws.Range("H46").Value = 1807.6923
$ws.Range("I46").Value = 925
$ws.Range("K46").Value = 925
$ws.Range("M46").Value = -737
$ws.Range("H82").Value = 5225.423
$ws.Range("I82").Value = 5575.2856
$ws.Range("K82").Value = 5575.2856
$ws.Range("M82").Value = -5214.2856
$ws.Range("H85").Value = 5225.423
$ws.Range("I85").Value = 5575.2856
$ws.Range("K85").Value = 5575.2856
$ws.Range("M85").Value = -4327.2856
$ws.Range("H122").Value = 5343.4287
$ws.Range("I122").Value = 2726
$ws.Range("J122").Value = 8833.333000000001
$ws.Range("K122").Value = 8178
$ws.Range("L122").Value = 26499.999
$ws.Range("M122").Value = -5728
$ws.Range("N122").Value = -31399.999
$ws.Range("H136").Value = 4120.6
$ws.Range("I136").Value = 1688.625
$ws.Range("J136").Value = 6900
$ws.Range("K136").Value = 5065.875
$ws.Range("L136").Value = 20700
$ws.Range("M136").Value = -2515.875
$ws.Range("N136").Value = -25800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 60762.363
$ws.Range("J46").Value = 60762.363
$ws.Range("L46").Value = 60762.363
$ws.Range("N46").Value = -61224.363
$ws.Range("H96").Value = 85319390
$ws.Range("I96").Value = 142858740
$ws.Range("J96").Value = 4764305
$ws.Range("K96").Value = 142858740
$ws.Range("L96").Value = 4764305
$ws.Range("M96").Value = -142857367
$ws.Range("N96").Value = -4767051
$ws.Range("H107").Value = 634.2727
$ws.Range("I107").Value = 560
$ws.Range("J107").Value = 764.25
$ws.Range("K107").Value = 1680
$ws.Range("L107").Value = 2292.75
$ws.Range("M107").Value = 240
$ws.Range("N107").Value = -6132.75
$ws.Range("H134").Value = 60762.363
$ws.Range("J134").Value = 60762.363
$ws.Range("L134").Value = 182287.089
$ws.Range("N134").Value = -187357.089
